$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 01:05"

# Update Estados Unidos (row 4) totals
$ws.Range("B4").Value = 1407396
$ws.Range("C4").Value = 21562
$ws.Range("D4").Value = 285976
$ws.Range("E4").Value = 1038109
$ws.Range("F4").Value = 16473
$ws.Range("G4").Value = 1516
$ws.Range("H4").Value = 83311

# Update El Salvador (row 100) totals
$ws.Range("E100").Value = 630
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 19

# Uruguay moves up in the ranking: row 116 now shows Uruguay with fresh data,
# Kenia drops to row 117 (keeping its prior data), Crucero drops to row 118
# (keeping its prior data).
$ws.Range("A116").Value = "Uruguay"
$ws.Range("B116").Value = 717
$ws.Range("C116").Value = 6
$ws.Range("D116").Value = 532
$ws.Range("E116").Value = 166
$ws.Range("F116").Value = 8
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 19

$ws.Range("A117").Value = "Kenia"
$ws.Range("B117").Value = 715
$ws.Range("C117").Value = 15
$ws.Range("D117").Value = 259
$ws.Range("E117").Value = 420
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 3
$ws.Range("H117").Value = 36

$ws.Range("A118").Value = "Crucero"
$ws.Range("B118").Value = 712
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 651
$ws.Range("E118").Value = 48
$ws.Range("F118").Value = 4
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 13
